$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add header "time_taken" in F1, copying the style used by the other headers (e.g. E1)
$ws.Range("F1").Value = "time_taken"
$ws.Range("E1").Copy()
$ws.Range("F1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# Fill F2:F25 with the recorded time_taken values as text strings
$timeTaken = @{
    2 = "2021-10-05 13:42:35.136607"
    3 = "2021-10-05 13:42:35.136620"
    4 = "2021-10-05 13:42:35.136624"
    5 = "2021-10-05 13:42:35.136627"
    6 = "2021-10-05 13:42:35.136630"
    7 = "2021-10-05 13:42:35.136634"
    8 = "2021-10-05 13:42:35.136636"
    9 = "2021-10-05 13:42:35.136639"
    10 = "2021-10-05 13:42:35.136643"
    11 = "2021-10-05 13:42:35.136646"
    12 = "2021-10-05 13:42:35.136649"
    13 = "2021-10-05 13:42:35.136652"
    14 = "2021-10-05 13:42:35.136655"
    15 = "2021-10-05 13:42:35.136658"
    16 = "2021-10-05 13:42:35.136660"
    17 = "2021-10-05 13:42:35.136663"
    18 = "2021-10-05 13:42:35.136667"
    19 = "2021-10-05 13:42:35.136670"
    20 = "2021-10-05 13:42:35.136673"
    21 = "2021-10-05 13:42:35.136676"
    22 = "2021-10-05 13:42:35.136679"
    23 = "2021-10-05 13:42:35.136682"
    24 = "2021-10-05 13:42:35.136685"
    25 = "2021-10-05 13:42:35.136688"
}

foreach ($row in $timeTaken.Keys) {
    $cell = $ws.Cells.Item($row, 6)
    $cell.Value = $timeTaken[$row]
}

